# Fix: ComboBoxEdit not updating after deleting a word
# Row 18 (item 14 - "Thuc hien...") had its actual-start/actual-finish
# dates (G18/H18) entered; apply the same date number-format already used
# by the other rows in those columns (e.g. G6/H6) and set the values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the date-format styling from an already-formatted cell (G6) and
# apply it to G18:H18 so they pick up the existing shared style instead of
# minting a new one.
$ws.Range("G6").Copy() | Out-Null
$ws.Range("G18:H18").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the actual start / actual finish dates (serials avoid re-triggering
# Excel's "looks like a date" autoformat, which would otherwise mint a
# duplicate number format).
$ws.Range("G18").Value = 43111   # 1/11/2018
$ws.Range("H18").Value = 43415   # 11/11/2018

# Move the active selection to H18 (was H15).
$ws.Range("H18").Select() | Out-Null
